$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 6: move the "X" mark from E6 (En proceso) to G6 (Terminado)
$ws.Range("E6").Value = $null
$ws.Range("G6").Value = "X"
$ws.Range("G6").HorizontalAlignment = -4152

# Row 7: move the "X" mark from E7 (En proceso) to G7 (Terminado)
$ws.Range("E7").Value = $null
$ws.Range("G7").Value = "X"
$ws.Range("G7").HorizontalAlignment = -4152

# Row 8: mark "En proceso" with an "X"
$ws.Range("E8").Value = "X"

# Row 9: mark "En proceso" with an "X"
$ws.Range("E9").Value = "X"

# Update the active selection to G8 (matches the saved view state)
$ws.Range("G8").Select() | Out-Null
